$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row for the extra "GetProductByName" endpoint in the
#    CATALOG API table (shifts old rows 6-16 down to 7-17, and the
#    "BASKET API" merged title row moves from A11:C11 to A12:C12
#    automatically).
# ---------------------------------------------------------------------------
$ws.Rows(6).Insert()

# Give the freshly inserted row the same plain bordered look the other
# data rows currently have, by copying formats from row 7 (an existing,
# still-untouched data row).
$ws.Range("A7:C7").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)

# Grab the same plain bordered look for the brand-new "ORDERING API" POST
# row (row 22), which keeps the old unfilled style.
$ws.Range("A7:C7").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Fill in cell values.
# ---------------------------------------------------------------------------

# --- CATALOG API table ---
$ws.Cells.Item(6,1).Value = "GET"
$ws.Cells.Item(6,2).Value = "api/v1/Catalog/GetProductByName/{productName}"
$ws.Cells.Item(6,3).Value = "Get matching products given the productName"

# --- ORDERING API (new section) ---
# Copy the "section title" format (font + grey fill + thin border, centred)
# from the BASKET API title row (row 12) onto the new title rows.
$ws.Range("A12:C12").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Range("A25:C25").PasteSpecial(-4122)

# Copy the "table header" format (font + grey fill + thin border, centred)
# from an existing "HTTP Method / Request URL / Use Case" row onto the new
# header rows.
$ws.Range("A13:C13").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Range("A26:C26").PasteSpecial(-4122)

$ws.Range("A19:C19").Merge()
$ws.Cells.Item(19,1).Value = "ORDERING API"

$ws.Cells.Item(20,1).Value = "HTTP Method"
$ws.Cells.Item(20,2).Value = "Request URL"
$ws.Cells.Item(20,3).Value = "Use Case"

$ws.Cells.Item(21,1).Value = "GET"
$ws.Cells.Item(21,2).Value = "api/v1/Order"
$ws.Cells.Item(21,3).Value = "Get orders for given username"

$ws.Cells.Item(22,1).Value = "POST"
$ws.Cells.Item(22,2).Value = "api/v1/Order"
$ws.Cells.Item(22,3).Value = "post order (just for testing-will not be used actually)"

# --- API GATEWAY MICROSERVICE (new section) ---
$ws.Range("A25:C25").Merge()
$ws.Cells.Item(25,1).Value = "API GATEWAY MICROSERVICE"

# Stamp a plain bordered look (same border used everywhere else) onto the
# new API GATEWAY MICROSERVICE data rows before colouring them below.
$ws.Range("A7:C7").Copy()
$ws.Range("A27:C33").PasteSpecial(-4122)

$ws.Cells.Item(26,1).Value = "HTTP Method"
$ws.Cells.Item(26,2).Value = "Request URL"
$ws.Cells.Item(26,3).Value = "Use Case"

$ws.Cells.Item(27,1).Value = "GET"
$ws.Cells.Item(27,2).Value = "/Catalog"
$ws.Cells.Item(27,3).Value = "api/v1/Catalog (Get)"

$ws.Cells.Item(28,1).Value = "POST"
$ws.Cells.Item(28,2).Value = "/Catalog"
$ws.Cells.Item(28,3).Value = "api/v1/Catalog (Post)"

$ws.Cells.Item(29,1).Value = "GET"
$ws.Cells.Item(29,2).Value = "/Catalog/{id}"
$ws.Cells.Item(29,3).Value = "api/v1/Catalog/{id}"

$ws.Cells.Item(30,1).Value = "GET"
$ws.Cells.Item(30,2).Value = "/Basket"
$ws.Cells.Item(30,3).Value = "api/v1/Basket (Get)"

$ws.Cells.Item(31,1).Value = "POST"
$ws.Cells.Item(31,2).Value = "/Basket"
$ws.Cells.Item(31,3).Value = "api/v1/Basket (Post)"

$ws.Cells.Item(32,1).Value = "POST"
$ws.Cells.Item(32,2).Value = "/Basket/Ckeckout"
$ws.Cells.Item(32,3).Value = "api/v1/Basket/Ckeckout (Post)"

$ws.Cells.Item(33,1).Value = "GET"
$ws.Cells.Item(33,2).Value = "/Order"
$ws.Cells.Item(33,3).Value = "api/v1/Order"

# ---------------------------------------------------------------------------
# 4. Re-colour every plain data row (currently unfilled, border-only) with
#    the new "accent1 / 40% lighter" blue fill used throughout the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A3:C9").Interior.ThemeColor = 5
$ws.Range("A14:C16").Interior.ThemeColor = 5
$ws.Range("A21:C21").Interior.ThemeColor = 5

# ---------------------------------------------------------------------------
# 5. Highlight the "checkout" row in yellow.
# ---------------------------------------------------------------------------
$ws.Range("A17:C17").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 6. Colour-code the new API GATEWAY MICROSERVICE rows: yellow for the
#    "pass-through" gateway routes, green for the two routes whose
#    downstream path differs from the incoming gateway path.
# ---------------------------------------------------------------------------
$ws.Range("A27:C28").Interior.Color = 65535
$ws.Range("A30:C31").Interior.Color = 65535
$ws.Range("A33:C33").Interior.Color = 65535

$ws.Range("A29:C29").Interior.Color = 5296274
$ws.Range("A32:C32").Interior.Color = 5296274

# ---------------------------------------------------------------------------
# 7. View tweaks: zoom to 130% and move the active selection to A17.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 130
$ws.Range("A17").Select()
